# The paragraph containing the literal text "<id>p026r_1</id>" was
# originally split across four separate runs ("<", "id>", "p026r_1",
# "</id>"). The edit merges them into a single run (keeping the
# formatting of the first of those runs: Courier New / color 7f6000 /
# sz 18) that holds the whole string "<id>p026r_1</id>".
#
# Doing a Find/Replace for the exact same text re-unifies the range
# Word selected (which spans the four original runs) into one run
# carrying the formatting of the range's first run - exactly matching
# the target edit.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "<id>p026r_1</id>",  # FindText
    $true,               # MatchCase
    $false,              # MatchWholeWord
    $false,              # MatchWildcards
    $false,              # MatchSoundsLike
    $false,              # MatchAllWordForms
    $true,                # Forward
    1,                   # Wrap (wdFindContinue)
    $false,              # Format
    "<id>p026r_1</id>",  # ReplaceWith (identical text -> merges runs)
    2                    # Replace (wdReplaceAll)
) | Out-Null
